$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.911561666666666
$ws.Range("H2").Value = 5.734684999999999
$ws.Range("I2").Value = 0.1720155802183755
$ws.Range("J2").Value = 0.1720155802183755
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.110264333333333
$ws.Range("N2").Value = 6.330793
$ws.Range("O2").Value = 0.3832041185227171
$ws.Range("P2").Value = 0.3832041185227171
$ws.Range("Q2").Value = 4.033900406133888
$ws.Range("R2").Value = 36.305103655205
$ws.Range("S2").Value = 0.06591707878975632
$ws.Range("T2").Value = 0.06591707878975629
$ws.Range("G3").Value = 1.911561666666666
$ws.Range("H3").Value = 5.734684999999999
$ws.Range("I3").Value = 0.1720155802183755
$ws.Range("J3").Value = 0.1720155802183755
$ws.Range("O3").Value = 0.0946183755984393
$ws.Range("P3").Value = 0.0946183755984393
$ws.Range("Q3").Value = 0.9960255782888886
$ws.Range("R3").Value = 8.964230204599998
$ws.Range("S3").Value = 0.01627583477788572
$ws.Range("T3").Value = 0.01627583477788571
$ws.Range("G4").Value = 1.911561666666666
$ws.Range("H4").Value = 5.734684999999999
$ws.Range("I4").Value = 0.1720155802183755
$ws.Range("J4").Value = 0.1720155802183755
$ws.Range("M4").Value = 2.065388333333333
$ws.Range("N4").Value = 6.196165
$ws.Range("O4").Value = 0.3750550597762889
$ws.Range("P4").Value = 0.3750550597762889
$ws.Range("Q4").Value = 3.948117164780555
$ws.Range("R4").Value = 35.53305448302499
$ws.Range("S4").Value = 0.06451531372125585
$ws.Range("T4").Value = 0.06451531372125582
$ws.Range("G5").Value = 1.911561666666666
$ws.Range("H5").Value = 5.734684999999999
$ws.Range("I5").Value = 0.1720155802183755
$ws.Range("J5").Value = 0.1720155802183755
$ws.Range("M5").Value = 0.8101876666666666
$ws.Range("N5").Value = 2.430563
$ws.Range("O5").Value = 0.1471224461025547
$ws.Range("P5").Value = 0.1471224461025547
$ws.Range("Q5").Value = 1.548723686406111
$ws.Range("R5").Value = 13.938513177655
$ws.Range("S5").Value = 0.02530735292947763
$ws.Range("T5").Value = 0.02530735292947763
$ws.Range("I6").Value = 0.4009917520372743
$ws.Range("J6").Value = 0.4009917520372743
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.110264333333333
$ws.Range("N6").Value = 6.330793
$ws.Range("O6").Value = 0.3832041185227171
$ws.Range("P6").Value = 0.3832041185227171
$ws.Range("Q6").Value = 9.403571405252888
$ws.Range("R6").Value = 84.632142647276
$ws.Range("S6").Value = 0.1536616908743237
$ws.Range("T6").Value = 0.1536616908743236
$ws.Range("I7").Value = 0.4009917520372743
$ws.Range("J7").Value = 0.4009917520372743
$ws.Range("O7").Value = 0.0946183755984393
$ws.Range("P7").Value = 0.0946183755984393
$ws.Range("S7").Value = 0.03794118820613906
$ws.Range("T7").Value = 0.03794118820613906
$ws.Range("I8").Value = 0.4009917520372743
$ws.Range("J8").Value = 0.4009917520372743
$ws.Range("M8").Value = 2.065388333333333
$ws.Range("N8").Value = 6.196165
$ws.Range("O8").Value = 0.3750550597762889
$ws.Range("P8").Value = 0.3750550597762889
$ws.Range("Q8").Value = 9.203598982975555
$ws.Range("R8").Value = 82.83239084678
$ws.Range("S8").Value = 0.1503939855301387
$ws.Range("T8").Value = 0.1503939855301387
$ws.Range("I9").Value = 0.4009917520372743
$ws.Range("J9").Value = 0.4009917520372743
$ws.Range("M9").Value = 0.8101876666666666
$ws.Range("N9").Value = 2.430563
$ws.Range("O9").Value = 0.1471224461025547
$ws.Range("P9").Value = 0.1471224461025547
$ws.Range("Q9").Value = 3.610285903435111
$ws.Range("R9").Value = 32.492573130916
$ws.Range("S9").Value = 0.05899488742667289
$ws.Range("T9").Value = 0.05899488742667287
$ws.Range("G10").Value = 4.603447666666667
$ws.Range("H10").Value = 13.810343
$ws.Range("I10").Value = 0.4142501574471451
$ws.Range("J10").Value = 0.4142501574471449
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.110264333333333
$ws.Range("N10").Value = 6.330793
$ws.Range("O10").Value = 0.3832041185227171
$ws.Range("P10").Value = 0.3832041185227171
$ws.Range("Q10").Value = 9.714491421333223
$ws.Range("R10").Value = 87.43042279199899
$ws.Range("S10").Value = 0.15874236643243
$ws.Range("T10").Value = 0.1587423664324299
$ws.Range("G11").Value = 4.603447666666667
$ws.Range("H11").Value = 13.810343
$ws.Range("I11").Value = 0.4142501574471451
$ws.Range("J11").Value = 0.4142501574471449
$ws.Range("O11").Value = 0.0946183755984393
$ws.Range("P11").Value = 0.0946183755984393
$ws.Range("Q11").Value = 2.398641751542222
$ws.Range("R11").Value = 21.58777576388
$ws.Range("S11").Value = 0.03919567698904659
$ws.Range("T11").Value = 0.03919567698904658
$ws.Range("G12").Value = 4.603447666666667
$ws.Range("H12").Value = 13.810343
$ws.Range("I12").Value = 0.4142501574471451
$ws.Range("J12").Value = 0.4142501574471449
$ws.Range("M12").Value = 2.065388333333333
$ws.Range("N12").Value = 6.196165
$ws.Range("O12").Value = 0.3750550597762889
$ws.Range("P12").Value = 0.3750550597762889
$ws.Range("Q12").Value = 9.507907103843889
$ws.Range("R12").Value = 85.57116393459499
$ws.Range("S12").Value = 0.1553666175636761
$ws.Range("T12").Value = 0.155366617563676
$ws.Range("G13").Value = 4.603447666666667
$ws.Range("H13").Value = 13.810343
$ws.Range("I13").Value = 0.4142501574471451
$ws.Range("J13").Value = 0.4142501574471449
$ws.Range("M13").Value = 0.8101876666666666
$ws.Range("N13").Value = 2.430563
$ws.Range("O13").Value = 0.1471224461025547
$ws.Range("P13").Value = 0.1471224461025547
$ws.Range("Q13").Value = 3.729656523678778
$ws.Range("R13").Value = 33.56690871310899
$ws.Range("S13").Value = 0.06094549646199241
$ws.Range("T13").Value = 0.06094549646199239
$ws.Range("G14").Value = 0.141604
$ws.Range("H14").Value = 0.424812
$ws.Range("I14").Value = 0.01274251029720526
$ws.Range("J14").Value = 0.01274251029720526
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.110264333333333
$ws.Range("N14").Value = 6.330793
$ws.Range("O14").Value = 0.3832041185227171
$ws.Range("P14").Value = 0.3832041185227171
$ws.Range("Q14").Value = 0.2988218706573333
$ws.Range("R14").Value = 2.689396835916
$ws.Range("S14").Value = 0.004882982426207188
$ws.Range("T14").Value = 0.004882982426207187
$ws.Range("G15").Value = 0.141604
$ws.Range("H15").Value = 0.424812
$ws.Range("I15").Value = 0.01274251029720526
$ws.Range("J15").Value = 0.01274251029720526
$ws.Range("O15").Value = 0.0946183755984393
$ws.Range("P15").Value = 0.0946183755984393
$ws.Range("Q15").Value = 0.07378323621333333
$ws.Range("R15").Value = 0.66404912592
$ws.Range("S15").Value = 0.001205675625367948
$ws.Range("T15").Value = 0.001205675625367948
$ws.Range("G16").Value = 0.141604
$ws.Range("H16").Value = 0.424812
$ws.Range("I16").Value = 0.01274251029720526
$ws.Range("J16").Value = 0.01274251029720526
$ws.Range("M16").Value = 2.065388333333333
$ws.Range("N16").Value = 6.196165
$ws.Range("O16").Value = 0.3750550597762889
$ws.Range("P16").Value = 0.3750550597762889
$ws.Range("Q16").Value = 0.2924672495533334
$ws.Range("R16").Value = 2.63220524598
$ws.Range("S16").Value = 0.004779142961218296
$ws.Range("T16").Value = 0.004779142961218295
$ws.Range("G17").Value = 0.141604
$ws.Range("H17").Value = 0.424812
$ws.Range("I17").Value = 0.01274251029720526
$ws.Range("J17").Value = 0.01274251029720526
$ws.Range("M17").Value = 0.8101876666666666
$ws.Range("N17").Value = 2.430563
$ws.Range("O17").Value = 0.1471224461025547
$ws.Range("P17").Value = 0.1471224461025547
$ws.Range("Q17").Value = 0.1147258143506667
$ws.Range("R17").Value = 1.032532329156
$ws.Range("S17").Value = 0.00187470928441183
$ws.Range("T17").Value = 0.001874709284411829
